$wb = $excel.ActiveWorkbook
$wb.Names.Add("zzz.tn.5250", "='#system'!`$Z`$2:`$Z`$6")
$n = $wb.Names.Item("zzz.tn.5250")
$n.Name = "tn.5250"
